# Practitioners have been parameterised.
# Adds a new "Practitioner" worksheet (parameterising practitioner test data,
# mirroring the existing "Organizations" parameter sheet), and updates the
# view state (active tab / zoom / selection) on the existing sheets.

$wb = $excel.ActiveWorkbook

$patients      = $wb.Worksheets.Item(1)
$organizations = $wb.Worksheets.Item(2)

# --- Add the new "Practitioner" sheet after the last existing sheet -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$practitioner = $wb.Worksheets.Add($null, $lastSheet)
$practitioner.Name = "Practitioner"

# Content: header + example code, blank row, then the "Must not exist" guard
# row (reusing the same shared string as the Organizations sheet) + value.
$practitioner.Range("A1").Value = "Practitioner Code"
$practitioner.Range("A2").Value = "G11111116"
$practitioner.Range("A4").Value = "Must not exist"
$practitioner.Range("A4").Font.Bold = $true
$practitioner.Range("A5").Value = "G99999999"

$practitioner.Range("A1").Select() | Out-Null

# --- View state updates ----------------------------------------------------

# Patients sheet: no longer the active tab; zoom 140 -> 120; scroll/selection.
$patients.Activate()
$excel.ActiveWindow.Zoom = 120
$patients.Range("B6").Select() | Out-Null

# Organizations sheet: zoom 140 -> 120; selection back to A1.
$organizations.Activate()
$excel.ActiveWindow.Zoom = 120
$organizations.Range("A1").Select() | Out-Null

# Practitioner sheet: active tab, zoom 120, selection at C26.
$practitioner.Activate()
$excel.ActiveWindow.Zoom = 120
$practitioner.Range("C26").Select() | Out-Null
